$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.313.42"
$ws.Range("E2").Value = "  +2.78%  "

$ws.Range("D3").Value = "3.406.02"

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.53%  "

$ws.Range("D8").Value = "3.395.61"
$ws.Range("E8").Value = "  +2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +13.10%  "

$ws.Range("E11").Value = "  +3.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.02"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.09%  "

$ws.Range("E13").Value = "  +6.40%  "

$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").Value = "3.948.87"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.41"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.415.36"
$ws.Range("E17").Value = "  +2.95%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("D19").Value = "65.260.01"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("E21").Value = "  +2.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "473.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +16.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +15.24%  "

$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.84%  "

$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.34"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.73"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.18%  "

$ws.Range("E32").Value = "  +2.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "575.31"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("E35").Value = "  +2.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.24%  "

$ws.Range("D39").Value = "0.0₃0764"
$ws.Range("E39").Value = "  +4.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.94"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.64%  "

$ws.Range("E41").Value = "  +2.42%  "

$ws.Range("D42").Value = "3.101.35"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("E44").Value = "  +2.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0417"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.86%  "

$ws.Range("E46").Value = "  +3.60%  "

$ws.Range("E47").Value = "  +6.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.08%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.28%  "

$ws.Range("E51").Value = "  +3.30%  "
